$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add attendance mark (1) for "săpt. 4" (column F) for the students who attended that week.
$rows = @(6, 9, 10, 12, 13, 14, 17, 18, 19, 21)
foreach ($r in $rows) {
    $ws.Cells.Item($r, 6).Value = 1
}

# Scroll the frozen pane so that the visible top-left unfrozen cell is C12 instead of C3.
$win = $excel.ActiveWindow
$win.ScrollRow = 12
$win.ScrollColumn = 3
